# Bitácora.xlsx update
# 1) Rename the existing sheet and add a new "Metricas" sheet after it.
# 2) Insert a new "Dificulty" column (C) on the activities sheet, fill it in,
#    and add a new activity row (row 12) with a couple of extra columns.
# 3) Add the "Metricas" sheet content (three small weekly tables).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename sheet / add Metricas sheet --------------------------------
$ws.Name = "Actividades y satisfacción"

$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Metricas"

# --- 2. Activities sheet: insert Dificulty column -------------------------
$ws.Columns("C").Insert()
$ws.Range("C4:C6").Merge()

$ws.Range("C4:C6").ClearFormats()
$ws.Range("C4:C6").HorizontalAlignment = -4131

$ws.Range("C2").Value = "Dificulty"
$ws.Range("C3").Value = "Easy"
$ws.Range("C4").Value = "Medium"
$ws.Range("C7").Value = "Medium"
$ws.Range("C8").Value = "Medium"
$ws.Range("C9").Value = "Easy"
$ws.Range("C10").Value = "Medium"
$ws.Range("C11").Value = "Medium"

# Column widths per the final layout (closest representable values - the
# host's column-width model snaps to the nearest 1/6 character unit)
$ws.Columns("C").ColumnWidth = 28.333333
$ws.Columns("D").ColumnWidth = 34.833333
$ws.Columns("E").ColumnWidth = 23.666667
$ws.Columns("F").ColumnWidth = 26.166667
$ws.Columns("G").ColumnWidth = 70.166667
$ws.Columns("H").ColumnWidth = 67.666667

# New row 12 (another Edwin / First delivery entry, with extra columns)
$r12 = $ws.Cells.Item(12, 1)
$r12.Value = 43523
$r12.NumberFormat = "d-mmm"
$ws.Cells.Item(12, 2).Value = "Edwin"
$ws.Cells.Item(12, 3).Value = "Easy"
$ws.Cells.Item(12, 4).Value = "First delivery"
$ws.Cells.Item(12, 5).Value = "Presentation"
$ws.Cells.Item(12, 6).Value = "30 mins (8:00 - 8:30)"
$ws.Cells.Item(12, 7).Value = "Nice, too easy, we had all the info"
$ws.Cells.Item(12, 8).Value = "Having the information is easier to make a presentation and takes less time"

$ws.Range("A2").Select()

# --- 3. Metricas sheet content --------------------------------------------
$ws2.Cells.Item(1, 1).Value = "Medicion de tiempos resultares semanales."

$ws2.Cells.Item(2, 1).Value = "Dificulty"
$ws2.Cells.Item(2, 2).Value = "Week 0"
$ws2.Cells.Item(2, 3).Value = "Week 1"
$ws2.Cells.Item(2, 4).Value = "Week 2"
$ws2.Cells.Item(2, 5).Value = "Week 3"
$ws2.Cells.Item(2, 6).Value = "Week 4"

$ws2.Cells.Item(3, 1).Value = "Easy"
$ws2.Cells.Item(3, 2).Value = "1 hour"

$ws2.Cells.Item(4, 1).Value = "Medium"
$ws2.Cells.Item(4, 2).Value = "2 hours"

$ws2.Cells.Item(5, 1).Value = "Hard"
$ws2.Cells.Item(5, 2).Value = "3 hours"

$ws2.Cells.Item(7, 1).Value = "Promedio de estado de ánimo"

$ws2.Cells.Item(8, 1).Value = "Person"
$ws2.Cells.Item(8, 2).Value = "Week 0"
$ws2.Cells.Item(8, 3).Value = "Week 1"
$ws2.Cells.Item(8, 4).Value = "Week 2"
$ws2.Cells.Item(8, 5).Value = "Week 3"
$ws2.Cells.Item(8, 6).Value = "Week 4"

$ws2.Cells.Item(9, 1).Value = "Edwin"
$ws2.Cells.Item(10, 1).Value = "Kirbey"
$ws2.Cells.Item(11, 1).Value = "Jorge"
$ws2.Cells.Item(12, 1).Value = "Mauricio"

$ws2.Cells.Item(14, 1).Value = "Contribuciónes del equipo"

$ws2.Cells.Item(15, 1).Value = "Person"
$ws2.Cells.Item(15, 2).Value = "Week 0"
$ws2.Cells.Item(15, 3).Value = "Week 1"
$ws2.Cells.Item(15, 4).Value = "Week 2"
$ws2.Cells.Item(15, 5).Value = "Week 3"
$ws2.Cells.Item(15, 6).Value = "Week 4"

$ws2.Cells.Item(16, 1).Value = "Edwin"
$ws2.Cells.Item(17, 1).Value = "Kirbey"
$ws2.Cells.Item(18, 1).Value = "Jorge"
$ws2.Cells.Item(19, 1).Value = "Mauricio"

$ws2.Cells.Item(21, 1).Value = "Velocidad del equipo"

$ws2.Cells.Item(22, 1).Value = "Person"
$ws2.Cells.Item(22, 2).Value = "Week 0"
$ws2.Cells.Item(22, 3).Value = "Week 1"
$ws2.Cells.Item(22, 4).Value = "Week 2"
$ws2.Cells.Item(22, 5).Value = "Week 3"
$ws2.Cells.Item(22, 6).Value = "Week 4"

$ws2.Cells.Item(23, 1).Value = "Edwin"
$ws2.Cells.Item(24, 1).Value = "Kirbey"
$ws2.Cells.Item(25, 1).Value = "Jorge"
$ws2.Cells.Item(26, 1).Value = "Mauricio"

$ws2.Range("C3").Select()

$ws.Activate()
